$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Los Ripudiatos"
$ws.Range("B10").Value = "Matteo Zanlucchi | SBARX"
$ws.Range("C10").Value = "Daniele Feltrinelli | Rita Levi’s"
$ws.Range("D10").Value = "Andrea Conzatti | FC SAVIGNANO"
$ws.Range("E10").Value = "Alessio Farinati | Pinguini Trentini"
$ws.Range("F10").Value = "Giovanni  Lasta | 4SINS"
